$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: "And can't rate because birth date, marital status are not
# brought in from json..." bullet gains "email," right after "because".
# ---------------------------------------------------------------------

# Re-write the span that used to be split around the "status are"
# grammar-check markers. Replacing it (even with identical text) makes
# the host collapse it back into a single run and drops the now-stale
# gramStart/gramEnd proofErr markers.
$d.Content.Find.Execute("marital status are not", $true, $false, $false, `
    $false, $false, $true, 1, $false, "marital status are not", 2) | Out-Null

# Insert " email," right after "because" in that same sentence.
$rngBecause = $d.Content
$rngBecause.Find.Execute("because", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rngBecause.Collapse(0)
$rngBecause.InsertAfter(" email,")

# ---------------------------------------------------------------------
# Part 2: four new bullet items after "Vehicle ownership, prior
# insurance ... during quote creation." and the trailing "_GoBack"
# bookmark ends up after the last of them.
# ---------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tail = $lastPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$newItems = @(
    "No way to just look up a quote",
    "Need valid driver’s license, state",
    "Comp and Collision required for Financed vehicle",
    "Not all errors returned are valid"
)

for ($i = 0; $i -lt $newItems.Length; $i++) {
    $para = $d.Paragraphs.Item($d.Paragraphs.Count)
    $r = $para.Range
    $r.Collapse(0)
    $r.InsertAfter($newItems[$i])
    if ($i -lt $newItems.Length - 1) {
        $para2 = $d.Paragraphs.Item($d.Paragraphs.Count)
        $r2 = $para2.Range
        $r2.Collapse(0)
        $r2.InsertParagraphAfter()
    }
}

# Re-create "_GoBack" right after the very last run. Adding a bookmark
# directly at a paragraph's text-end offset mis-resolves in this host,
# so park a throw-away character there first, bookmark across it (a
# non-degenerate range resolves correctly), then delete the character --
# the bookmark collapses back down to the correct zero-width spot.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$finalRange.Collapse(0)
$finalRange.InsertAfter("~")
$markEnd = $finalRange.End
$markStart = $markEnd - 1
$safeRange = $d.Range($markStart, $markEnd)
$d.Bookmarks.Add("_GoBack", $safeRange) | Out-Null
$d.Range($markStart, $markEnd).Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
